$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "50.614.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.909.92"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "373.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.15"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.04%  "
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.54"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0838"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.372.06"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.47"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.909.18"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.00"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +49.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.978"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "50.603.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("E20").Value = "  -6.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0945"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.57"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.29"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("E25").Value = "  +7.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.93"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.27"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.25%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.30"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("E31").Value = "  -8.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.89"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.67"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0432"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.76%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.04"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.28"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.77"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("E42").Value = "  -5.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.78"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.69"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.04"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.971.49"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0321"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.14"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.91%  "
